$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows 10:11 -- Lower/Upper 95% HDI summary rows appended below the
# existing mean/sd/CV/quantile block (rows 1-9 are left untouched)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Lower 95% HDI"
$ws.Range("A10").Font.Bold = $true

$ws.Range("B10:C11").NumberFormat = "#,##0"

$ws.Range("B10").Value = 110000
$ws.Range("C10").Value = 172000

$ws.Range("A11").Value = "Upper 95% HDI"
$ws.Range("A11").Font.Bold = $true

$ws.Range("B11").Value = 328000
$ws.Range("C11").Value = 256000

# ---------------------------------------------------------------------------
# Row 17 -- header row for the wide summary table (raw values)
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "mean"
$ws.Range("C17").Value = "sd"
$ws.Range("D17").Value = "CV"
$ws.Range("E17").Value = 0.025
$ws.Range("F17").Value = 0.25
$ws.Range("G17").Value = 0.5
$ws.Range("H17").Value = 0.75
$ws.Range("I17").Value = 0.975
$ws.Range("J17").Value = "Lower 95% HDI"
$ws.Range("K17").Value = "Upper 95% HDI"

# Row 18 -- AR1-Empirical
$ws.Range("A18").Value = "AR1-Empirical"
$ws.Range("B18").Value = 213317.467815343
$ws.Range("C18").Value = 58060.452444853399
$ws.Range("D18").Value = 0.27217861265404242
$ws.Range("E18").Value = 121446.57067827
$ws.Range("F18").Value = 171732.38381670101
$ws.Range("G18").Value = 205943.77193442601
$ws.Range("H18").Value = 246675.25088838401
$ws.Range("I18").Value = 347561.76785331802
$ws.Range("J18").Value = 110000
$ws.Range("K18").Value = 328000

# Row 19 -- AR1-Base
$ws.Range("A19").Value = "AR1-Base"
$ws.Range("B19").Value = 213361.25441453099
$ws.Range("C19").Value = 21546.309475039401
$ws.Range("D19").Value = 0.10098510872633859
$ws.Range("E19").Value = 174124.38220897
$ws.Range("F19").Value = 198274.72760430499
$ws.Range("G19").Value = 212369.381770754
$ws.Range("H19").Value = 227248.057875008
$ws.Range("I19").Value = 258642.84210630201
$ws.Range("J19").Value = 172000
$ws.Range("K19").Value = 256000

# ---------------------------------------------------------------------------
# Row 24 -- header row for the condensed (/1000) summary table
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "Model"
$ws.Range("B24").Value = "mean"
$ws.Range("C24").Value = "sd"
$ws.Range("D24").Value = "CV"
$ws.Range("E24").Value = 0.025
$ws.Range("F24").Value = 0.25
$ws.Range("G24").Value = 0.5
$ws.Range("H24").Value = 0.75
$ws.Range("I24").Value = 0.975
$ws.Range("J24").Value = "Lower 95% HDI"
$ws.Range("K24").Value = "Upper 95% HDI"

# Row 25/26 values divided by 1000, formatted with numFmtId 1 ("0")
$ws.Range("A25").Value = "AR1-Empirical"
$ws.Range("A26").Value = "AR1-Base"

$ws.Range("B25:K26").NumberFormat = "0"

$ws.Range("B25:B26").FormulaR1C1 = "=R[-7]C/1000"
$ws.Range("C25:K26").FormulaR1C1 = "=R[-7]C/1000"

# ---------------------------------------------------------------------------
# Cosmetic sheet adjustments
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 13

$ws.Range("A24:K26").Select() | Out-Null
